$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 92
    3  = 86
    4  = 100
    5  = 112
    6  = 106
    7  = 114
    8  = 86
    9  = 72
    10 = 118
    11 = 116
    12 = 106
    13 = 92
    14 = 100
    15 = 78
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 9).Value = $values[$row]
}

$ws.Range("I15").Select()
